# Append the latest profit row (2025-11-20 run) to the bottom of the
# Date/Profit table on the single worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95

# Write the date as literal text (matching every other date cell in the
# column, which are stored as plain strings, not date serials). Forcing
# NumberFormat to "@" (Text) before the assignment stops Excel's normal
# autodetection from turning "11/20/2025" into a date serial number; the
# ClearFormats() afterwards removes the now-unneeded text format override so
# the cell is left with the default (unstyled) formatting, same as its
# neighbors.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "11/20/2025"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value = 8188.95
